$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Revision-history table: "Elaboration iteration 3 draft 1" ->
#    "Construction iteration 4 draft 1" (Construction-phase iteration bump).
# ---------------------------------------------------------------------------
$range = $d.Content
$found = $range.Find.Execute("Elaboration iteration 3 draft 1", $false, $false,
                              $false, $false, $false, $true, 1, $false,
                              "Construction iteration 4 draft 1", 2)

# ---------------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark (Word's "last edit location" marker) from
#    the old edit spot (inside "Hyppighed", next to "kundeoplysninger") to
#    the new edit spot, right after "(Scope)" in the "Afgraensning (Scope)"
#    heading - this reflects where the author's final edit of the session
#    landed.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$scopeRange = $d.Content
$scopeFound = $scopeRange.Find.Execute("(Scope)")
if ($scopeFound) {
    $scopeRange.Collapse(0)
    $d.Bookmarks.Add("_GoBack", $scopeRange)
}
